$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for each data row (rows 2-43).
# Update every one of those cells by incrementing the stored date by one day
# (serial 45801 -> 45802, i.e. 2025-05-24 -> 2025-05-25).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -ne $null) {
        $cell.Value = $cell.Value2 + 1
    }
}
